# Update gh-pages to output generated at 456a3b4
# Updates the "想去人数" (F column) counts on the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 2087
$ws1.Range("F6").Value  = 637
$ws1.Range("F7").Value  = 103
$ws1.Range("F8").Value  = 2075
$ws1.Range("F9").Value  = 10717
$ws1.Range("F14").Value = 421
$ws1.Range("F15").Value = 7574
$ws1.Range("F18").Value = 268
$ws1.Range("F19").Value = 69
$ws1.Range("F20").Value = 3341

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 2087
$ws4.Range("F6").Value  = 637
$ws4.Range("F8").Value  = 103
$ws4.Range("F9").Value  = 2075
$ws4.Range("F12").Value = 10717
$ws4.Range("F17").Value = 421
$ws4.Range("F18").Value = 7574
$ws4.Range("F21").Value = 268
$ws4.Range("F22").Value = 69
$ws4.Range("F23").Value = 3341
